$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.62612783908844
$ws.Range("B1").Value = 2.142900943756104
$ws.Range("C1").Value = 2.315966129302979
$ws.Range("D1").Value = 2.965813636779785
$ws.Range("E1").Value = 2.083703279495239
